# "add see doc and upload" — duplicate the Use-Case table (cols A:E) into a
# second copy anchored at column H, dropping the "User" column (D) so the
# new block only carries No. / Use Case / Admin / Keterangan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- copy each source column into its destination column (values + styles) ---
$ws.Range("A1:A15").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4104) | Out-Null   # xlPasteAll
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B1:B15").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4104) | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null

$ws.Range("C1:C15").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4104) | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null

$ws.Range("E1:E15").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4104) | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- row heights grew (taller wrapped text / merged header) after the paste ---
$ws.Rows(1).RowHeight = 28.8
$ws.Rows(2).RowHeight = 43.2
$ws.Rows(3).RowHeight = 72
$ws.Rows(4).RowHeight = 57.6
$ws.Rows(5).RowHeight = 72
$ws.Rows(6).RowHeight = 100.8
$ws.Rows(7).RowHeight = 86.4
$ws.Rows(8).RowHeight = 86.4
$ws.Rows(9).RowHeight = 72
$ws.Rows(10).RowHeight = 129.6
$ws.Rows(11).RowHeight = 86.4
$ws.Rows(12).RowHeight = 72
$ws.Rows(13).RowHeight = 100.8
$ws.Rows(14).RowHeight = 129.6
$ws.Rows(15).RowHeight = 86.4

# --- final selection matches the new pasted block ---
$ws.Range("H1:K15").Select() | Out-Null
